$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 14888
$ws.Range("F3").Value = 18566
$ws.Range("F14").Value = 110
$ws.Range("F21").Value = 229
$ws.Range("F22").Value = 7697
$ws.Range("F26").Value = 1222
$ws.Range("F28").Value = 5961
$ws.Range("F29").Value = 103
$ws.Range("F30").Value = 64
$ws.Range("F32").Value = 154
$ws.Range("F34").Value = 5317
$ws.Range("F36").Value = 40
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 2
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 14888
$ws.Range("F3").Value = 18566
$ws.Range("F14").Value = 110
$ws.Range("F22").Value = 229
$ws.Range("F23").Value = 7697
$ws.Range("F27").Value = 1222
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 5961
$ws.Range("F32").Value = 103
$ws.Range("F33").Value = 64
$ws.Range("F35").Value = 154
$ws.Range("F37").Value = 5317
$ws.Range("F39").Value = 40

$wb.Save()
